$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = '$patient: Patient'
$ws.Range("B8").Select()
